$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) DEMOS section: note text updates (no row/structure changes here)
# ---------------------------------------------------------------------
$ws.Range("D9").Value  = "Point unique alors qu'on devrait voir une ligne"
$ws.Range("D13").Value = "Le nouveau dispatch semble avoir ralenti considérablement cette ROM"
$ws.Range("D15").Value = "N/A"
$ws.Range("D16").Value = "Le nouveau dispatch semble avoir ralenti considérablement cette ROM"

# ---------------------------------------------------------------------
# 2) GAMES section restructure
#    Old: alien, Herdle, MusicMaker, MegaManX16, Ninja, Pong, Reflection, Snafu, tetris
#    New: alien, Herdle, MusicMaker, Ninja, Pacman, Pong, Reflection, Snafu
#    A single row delete (MegaManX16) plus a few cell overwrites reproduces
#    the whole restructure without ever creating new style entries.
# ---------------------------------------------------------------------
$ws.Rows(26).Delete()

# MusicMaker note
$ws.Range("D25").Value = "Jeu de son"

# Ninja: status goes back to "fine" (green) and note is cleared
$ws.Range("C5").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D26").Value = "N/A"

# Shift the remaining names down one slot, inserting Pacman
$ws.Range("B27").Value = "Pacman"
$ws.Range("B28").Value = "Pong"
$ws.Range("B29").Value = "Reflection"
$ws.Range("B30").Value = "Snafu"

# Reflection row (now holding former Snafu data) & Snafu row (former tetris data)
# need their status recolored back to green and notes cleared
$ws.Range("C5").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D29").Value = "N/A"

$ws.Range("C5").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("D30").Value = "N/A"

# ---------------------------------------------------------------------
# 3) Cosmetic: column widths & selection
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 1.6666666666666667
$ws.Columns("D").ColumnWidth = 67.66666666666667
$ws.Columns("E").ColumnWidth = 1

$ws.Range("F2").Select()
